# Partial check-in during mind-numbing edit. push survey Uri construction into
# Java side so it is implemented just once.  Add odkCommon.closeWindow API.
#
# This updates the "queries" sheet: the old "auxillaryHash" column (G) is
# renamed to "newRowInitialElementKeyToValueMap" and its JS-ish URI-building
# value is replaced by a simple object-literal expression; a new
# "openRowInitialElementKeyToValueMap" column (H) is added with value "{}".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("queries")

# Make "queries" the active sheet/tab (this also clears tabSelected on
# whichever sheet was previously active, and updates workbook activeTab).
$ws.Activate()

# Widen the new column H to roughly match the source column width.
$ws.Columns.Item(8).ColumnWidth = 33.6

# Row heights for the header/data rows.
$ws.Rows.Item(1).RowHeight = 19.65
$ws.Rows.Item(2).RowHeight = 25.55

# Header row: add header for new column H first, then rename column G, so
# that new shared-string entries are appended in the same order as the
# source workbook.
$ws.Range("H1").Value2 = "openRowInitialElementKeyToValueMap"
$ws.Range("G1").Value2 = "newRowInitialElementKeyToValueMap"

# Data row: add the new H2 value "{}" first, then replace the old ad-hoc
# URI-construction snippet in G2 with the simple object literal, preserving
# the cell's quote-prefix text style by prefixing with a literal leading
# apostrophe (Excel's "treat as text" marker, which is not stored as part
# of the text itself).
$ws.Range("H2").Value2 = "{}"
$ws.Range("G2").Value = "'{ household_id : data('household_id') }"

# Put the selection where the source file shows it afterwards.
$ws.Range("D8").Select()
